# Generate Report for Handoff
# This script updates the localization-status workbook to reflect that
# "b.md" is now "Ready for handoff" (instead of "Handed back: in sync
# with en-US"), with a fresh handoff timestamp and a new handoff
# artifact file name (b.63290e5768f688058c7b37413b0a5c26c308f864.*).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" row (B=zh-cn status, C=de-de
# status, D=Latest Handoff Date)
# ---------------------------------------------------------------------
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-27-21 06:27:53"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" row
# ---------------------------------------------------------------------
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-21 06:27:49"

# Rebuild the hyperlinks on the zh-cn sheet so the "display" text for D3
# matches the new file name while every other hyperlink (and the
# underlying link targets) stay exactly as they were. (Deleting the
# hyperlinks anywhere on the sheet clears all of them, so we recreate
# the full set in the original order/addresses.)
$wsZhCn.Range("A1").Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/211f9b0f23ea1c7030d69b52d34ea59beafbe0bf/e2e/a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/211f9b0f23ea1c7030d69b52d34ea59beafbe0bf/e2e/a.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aeb18ce7bd2a048ed0442fea30f75bb5c8f436a1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3e449d305f28b6af11f153c3b191f5dd735c7f13/e2e/a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a9d32ec142c962b949a27aae4fe509f25249652c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/211f9b0f23ea1c7030d69b52d34ea59beafbe0bf/e2e/b.md", "", "", "b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/211f9b0f23ea1c7030d69b52d34ea59beafbe0bf/e2e/b.md", "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aeb18ce7bd2a048ed0442fea30f75bb5c8f436a1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3e449d305f28b6af11f153c3b191f5dd735c7f13/e2e/a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a9d32ec142c962b949a27aae4fe509f25249652c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" row
# ---------------------------------------------------------------------
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-21 06:27:53"

$wsDeDe.Range("A1").Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/211f9b0f23ea1c7030d69b52d34ea59beafbe0bf/e2e/a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/211f9b0f23ea1c7030d69b52d34ea59beafbe0bf/e2e/a.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a81eb5ec4d8904ff3de196cb1a73563392c1bde8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6f1648ef52841e3ca690c7d357ac99621dd017fa/e2e/a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/429d7720526c1292c10e770f78c674142b94609a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/211f9b0f23ea1c7030d69b52d34ea59beafbe0bf/e2e/b.md", "", "", "b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/211f9b0f23ea1c7030d69b52d34ea59beafbe0bf/e2e/b.md", "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a81eb5ec4d8904ff3de196cb1a73563392c1bde8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6f1648ef52841e3ca690c7d357ac99621dd017fa/e2e/a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/429d7720526c1292c10e770f78c674142b94609a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")

Write-Host "Report regenerated for handoff."
